$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old trailing row (content formerly in row 24) - the sheet shrinks
# from A1:C24 to A1:C23.
$ws.Rows.Item(24).Delete()

# Row 13: "Programa resumido:" / "Semestral"
$ws.Cells.Item(13,1).Value2 = 'Programa resumido:'
$ws.Cells.Item(13,2).Value2 = 'Semestral'
$ws.Cells.Item(13,3).Value2 = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

# Row 14: "Short syllabus:" only
$ws.Cells.Item(14,1).Value2 = 'Short syllabus:'
$ws.Cells.Item(14,2).ClearContents()
$ws.Cells.Item(14,3).ClearContents()
$ws.Rows.Item(14).RowHeight = 60

# Row 15: "Programa:" / "01/01/2012"
$ws.Cells.Item(15,1).Value2 = 'Programa:'
$ws.Cells.Item(15,2).Value2 = '01/01/2012'
$ws.Cells.Item(15,3).Value2 = '01/01/2012'
$ws.Rows.Item(15).RowHeight = 120

# Row 16: "Syllabus:" only
$ws.Cells.Item(16,1).Value2 = 'Syllabus:'
$ws.Cells.Item(16,2).ClearContents()
$ws.Cells.Item(16,3).ClearContents()
$ws.Rows.Item(16).RowHeight = 120

# Row 17: "Avaliação:" only
$ws.Cells.Item(17,1).Value2 = 'Avaliação:'
$ws.Cells.Item(17,2).ClearContents()
$ws.Cells.Item(17,3).ClearContents()
$ws.Rows.Item(17).RowHeight = 15

# Row 18: "Método:" / "519033 - Carlos Yujiro Shigue"
$ws.Cells.Item(18,1).Value2 = 'Método:'
$ws.Cells.Item(18,2).Value2 = '519033 - Carlos Yujiro Shigue'
$ws.Cells.Item(18,3).Value2 = '519033 - Carlos Yujiro Shigue'
$ws.Rows.Item(18).RowHeight = 60

# Row 19: "Critério:" / evaluation method text
$ws.Cells.Item(19,1).Value2 = 'Critério:'
$ws.Cells.Item(19,2).Value2 = 'A avaliação será feita por meio de duas provas escritas P1 e P2 e por listas de exercícios e relatórios.'
$ws.Cells.Item(19,3).Value2 = 'A avaliação será feita por meio de duas provas escritas P1 e P2 e por listas de exercícios e relatórios.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20: "Norma de recuperação:" / final grade formula
$ws.Cells.Item(20,1).Value2 = 'Norma de recuperação:'
$ws.Cells.Item(20,2).Value2 = 'A Nota final (NF) será calculada pela média ponderada das provas escritas e pela média dos trabalhos TR da seguinte maneira: NF = (P1 + 2*P2 + TR)/4'
$ws.Cells.Item(20,3).Value2 = 'A Nota final (NF) será calculada pela média ponderada das provas escritas e pela média dos trabalhos TR da seguinte maneira: NF = (P1 + 2*P2 + TR)/4'
$ws.Rows.Item(20).RowHeight = 120

# Row 21: "Bibliografia:" / makeup exam formula
$ws.Cells.Item(21,1).Value2 = 'Bibliografia:'
$ws.Cells.Item(21,2).Value2 = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Cells.Item(21,3).Value2 = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Rows.Item(21).RowHeight = 15

# Row 22: "Requisitos:" only
$ws.Cells.Item(22,1).Value2 = 'Requisitos:'
$ws.Cells.Item(22,2).ClearContents()
$ws.Cells.Item(22,3).ClearContents()
$ws.Rows.Item(22).RowHeight = 15

# Row 23: prerequisite text (B/C only, no A label)
$ws.Cells.Item(23,1).ClearContents()
$ws.Cells.Item(23,2).Value2 = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Cells.Item(23,3).Value2 = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30
